$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '66.873.24'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Range('E2').Value = '  +6.08%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.539.63'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Range('E3').Value = '  +9.59%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '568.84'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Range('E5').Value = '  +7.83%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '189.13'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Range('E6').Value = '  +10.92%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.534.33'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Range('E7').Value = '  +9.44%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.620'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Range('E8').Value = '  +4.07%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.635'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Range('E10').Value = '  +5.08%  '
$ws.Range('E11').Value = '  +13.43%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '54.89'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Range('E12').Value = '  +3.81%  '
$ws.Range('E13').Value = '  +7.11%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '9.43'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Range('E14').Value = '  +3.22%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.102.50'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Range('E15').Value = '  +9.68%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.537.35'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Range('E16').Value = '  +9.69%  '
$ws.Range('E17').Value = '  +4.66%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '66.907.76'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Range('E18').Value = '  +6.45%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '18.26'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Range('E19').Value = '  +6.32%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '12.05'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Range('E20').Value = '  +8.89%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.999'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Range('E21').Value = '  +3.32%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '430.83'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Range('E22').Value = '  +17.65%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '4.22'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Range('E23').Value = '  +12.95%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '85.17'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Range('E24').Value = '  +5.07%  '
$ws.Range('E25').Value = '  +4.08%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '11.12'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.91'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Range('E27').Value = '  +10.42%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '12.29'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Range('E28').Value = '  +9.32%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '9.26'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Range('E29').Value = '  +13.03%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '30.39'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Range('E30').Value = '  +6.89%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '642.37'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '6.63'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Range('E32').Value = '  +3.19%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '11.77'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Range('E33').Value = '  +5.28%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.113'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Range('E34').Value = '  +6.99%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '59.88'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Range('E35').Value = '  +5.60%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.151'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Range('E36').Value = '  +23.18%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '38.71'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Range('E37').Value = '  +5.71%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.0₃0817'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Range('E38').Value = '  +14.79%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.392'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Range('E40').Value = '  +4.11%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.35'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Range('E41').Value = '  +14.32%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.00'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '3.055.95'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Range('E43').Value = '  +6.58%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.68'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Range('E44').Value = '  +4.99%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.88'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Range('E45').Value = '  +11.91%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.34'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Range('E46').Value = '  +8.72%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0421'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Range('E47').Value = '  +7.54%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.79'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Range('E48').Value = '  +4.48%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.132'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Range('E49').Value = '  +5.77%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '141.46'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Range('E50').Value = '  +5.67%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '8.65'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Range('E51').Value = '  +11.16%  '
